$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at row 552, pushing the existing rows 552-588
# down to 554-590 (weekly update: two new price observations added).
$ws.Rows.Item(552).Insert()
$ws.Rows.Item(552).Insert()

# New row 552
$ws.Cells.Item(552,1).Value = 11
$ws.Cells.Item(552,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(552,3).Value = "Bíobío"
$ws.Cells.Item(552,4).Value = 44931
$ws.Cells.Item(552,5).Value = 8
$ws.Cells.Item(552,6).Value = 100112020
$ws.Cells.Item(552,7).Value = "Tomate"
$ws.Cells.Item(552,8).Value = "Larga vida"
$ws.Cells.Item(552,9).Value = "Primera"
$ws.Cells.Item(552,10).Value = 270
$ws.Cells.Item(552,11).Value = 14000
$ws.Cells.Item(552,12).Value = 15000
$ws.Cells.Item(552,13).Value = 14444
$ws.Cells.Item(552,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(552,15).Value = "Provincia de Quillota"
$ws.Cells.Item(552,16).Value = 802
$ws.Cells.Item(552,17).Value = 18
$ws.Cells.Item(552,18).Value = "Hortaliza"

# New row 553
$ws.Cells.Item(553,1).Value = 11
$ws.Cells.Item(553,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(553,3).Value = "Bíobío"
$ws.Cells.Item(553,4).Value = 44931
$ws.Cells.Item(553,5).Value = 8
$ws.Cells.Item(553,6).Value = 100112020
$ws.Cells.Item(553,7).Value = "Tomate"
$ws.Cells.Item(553,8).Value = "Larga vida"
$ws.Cells.Item(553,9).Value = "Primera"
$ws.Cells.Item(553,10).Value = 450
$ws.Cells.Item(553,11).Value = 4500
$ws.Cells.Item(553,12).Value = 5000
$ws.Cells.Item(553,13).Value = 4778
$ws.Cells.Item(553,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(553,15).Value = "Quillón"
$ws.Cells.Item(553,16).Value = 478
$ws.Cells.Item(553,17).Value = 10
$ws.Cells.Item(553,18).Value = "Hortaliza"
